# "Aula 10 - Python Big Data" — update the title placeholder on slide 2:
#   "Aula 01"          -> "Aula 10"
#   "Contextualização" -> "Projeto AWS " + "Big Data" (two separate runs)

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(2)
$sh = $s.Shapes.Item(3)          # "Título 2" placeholder shape
$tr = $sh.TextFrame.TextRange

# --- "Aula 01" -> "Aula 10" --------------------------------------------
$aula = $tr.Find("Aula 01", 0)
$aula.Text = "Aula 10"

# --- "Contextualização" -> "Projeto AWS " + "Big Data" -----------------
$newText     = "Projeto AWS Big Data"
$bigDataText = "Big Data"

$ctx      = $tr.Find("Contextualização", 0)
$ctxStart = $ctx.Start
$ctx.Text = $newText

# Re-apply formatting to just the trailing "Big Data" portion so it is
# written out as its own run, distinct from the leading "Projeto AWS ".
$bigDataStart = $ctxStart + $newText.Length - $bigDataText.Length
$bigData = $tr.Characters($bigDataStart, $bigDataText.Length)
$bigData.Font.Size = 36
$bigData.Font.Bold = $true
